# Rename sheets and update summary values to create a rounded "compare tab" view.
$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wb.Worksheets.Item(1).Name = "data_CCM code_FLASH"
$wb.Worksheets.Item(2).Name = "data_CCM code_CCM"
$wb.Worksheets.Item(3).Name = "data_RAM code_FLASH"
$wb.Worksheets.Item(4).Name = "data_RAM code_CCM"

# --- Sheet 1: data_CCM code_FLASH ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = 13658
$ws1.Range("C2").Value = 23686
$ws1.Range("D2").Value = 31501
$ws1.Range("B5").Value = 42.575
$ws1.Range("C5").Value = 43.107
$ws1.Range("D5").Value = 45.107

# --- Sheet 2: data_CCM code_CCM ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = 11852
$ws2.Range("C2").Value = 22982
$ws2.Range("D2").Value = 34025
$ws2.Range("B5").Value = 36.953
$ws2.Range("C5").Value = 35.847
$ws2.Range("D5").Value = 35.398

# --- Sheet 3: data_RAM code_FLASH ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = 14037
$ws3.Range("C2").Value = 24263
$ws3.Range("D2").Value = 31891
$ws3.Range("B5").Value = 42.077
$ws3.Range("C5").Value = 42.742
$ws3.Range("D5").Value = 44.416

# --- Sheet 4: data_RAM code_CCM ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = 12113
$ws4.Range("C2").Value = 23487
$ws4.Range("D2").Value = 34846
$ws4.Range("B5").Value = 36.319
$ws4.Range("C5").Value = 35.19
$ws4.Range("D5").Value = 34.829
